$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (Coin name / Link / Volume columns) -- not numeric-looking,
# so a normal .Value assignment keeps them as text.
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E9").Value = "8OneONEBestin24h"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"

# Numeric-looking text updates (Price / Hora columns). These must stay plain text
# (the sheet stores them as inline strings, e.g. "0.005020" with the trailing zero,
# "--" placeholders elsewhere, etc.), so force text format before assigning, then
# drop back to the default style (no NumberFormat override survives on the cell).
$textCells = @{
    "D2" = "247.55"
    "G2" = "7"
    "D3" = "21.75"
    "G3" = "7"
    "D4" = "5.458"
    "G4" = "7"
    "D5" = "0.05665"
    "G5" = "7"
    "D6" = "3.372"
    "G6" = "7"
    "D7" = "0.8006"
    "G7" = "7"
    "G8" = "7"
    "D9" = "0.01155"
    "G9" = "7"
    "D10" = "0.1448"
    "G10" = "7"
    "D11" = "0.07226"
    "G11" = "7"
    "D12" = "0.03148"
    "G12" = "7"
    "D13" = "0.02936"
    "G13" = "7"
    "D14" = "0.09284"
    "G14" = "7"
    "D15" = "0.001639"
    "G15" = "7"
    "D16" = "3.226"
    "G16" = "7"
    "D17" = "0.04726"
    "G17" = "7"
    "D18" = "0.006447"
    "G18" = "7"
    "D19" = "0.005020"
    "G19" = "7"
    "G20" = "7"
    "G21" = "7"
    "G22" = "7"
    "D23" = "3.830"
    "G23" = "7"
    "D24" = "6.435"
    "G24" = "7"
    "D25" = "2.086"
    "G25" = "7"
    "D26" = "0.3277"
    "G26" = "7"
    "G27" = "7"
    "G28" = "7"
    "G29" = "7"
    "G30" = "7"
    "G31" = "7"
    "G32" = "7"
    "G33" = "7"
    "G34" = "7"
    "G35" = "7"
    "G36" = "7"
    "G37" = "7"
    "G38" = "7"
    "G39" = "7"
    "D40" = "0.04090"
    "G40" = "7"
    "D41" = "0.006930"
    "G41" = "7"
    "D42" = "0.1040"
    "G42" = "7"
    "D43" = "0.003300"
    "G43" = "7"
    "D44" = "0.009032"
    "G44" = "7"
    "G45" = "7"
    "G46" = "7"
    "G47" = "7"
    "D48" = "0.01055"
    "G48" = "7"
    "G49" = "7"
    "G50" = "7"
    "G51" = "7"
}
foreach ($ref in $textCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$ref]
    $cell.Style = "Normal"
}
